$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 12 (old row 12 "Programa resumido:" and everything
# below it shifts down by 4 rows, to rows 16-26)
$ws.Rows("12:15").Insert()

# The insert operation fills the new rows with column A's formatting (copied down
# from row 11). Clear that stray formatting in column A for rows 13-15, since in
# the target those rows have no value/format in column A at all.
$ws.Range("A13:A15").Clear()

# New row 12: header only in column A
$ws.Range("A12").Value = "Docentes responsáveis:"

# New rows 13-15: docentes names in both B and C (no A, no custom row height)
$ws.Range("B13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

$ws.Range("B14").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C14").Value = "2166002 - Sandra Giacomin Schneider"

$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"

# Fix the formats on B13:C15 to match the normal column B/C wrap-text styles used
# throughout the sheet (reuse the existing styles rather than Excel auto-picking a
# bold header style for them).
$ws.Range("B10").Copy()
$ws.Range("B13:B15").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13:C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
